$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a few odds values in row 2
$ws.Range("G2").Value = 3.3
$ws.Range("I2").Value = 2.55
$ws.Range("X2").Value = 15
$ws.Range("AH2").Value = 10
$ws.Range("AJ2").Value = 26

# Remove row 4 (Sunderland - West Brom) entirely; row 5 (Wanderers - Miramar) shifts up to row 4
$ws.Rows(4).Delete()
